{"js": "// Replace the 100 \"NN\u00d7NN=\" multiplication prompts in the single 20x5 table\n// with the new values from the target revision. Cell formatting (fonts,\n// size, paragraph alignment) is left untouched because we only write the\n// text value of each cell, not its structure.\nconst newValues = [\n  [\"97\u00d724=\", \"27\u00d774=\", \"34\u00d785=\", \"61\u00d735=\", \"49\u00d721=\"],\n  [\"54\u00d758=\", \"43\u00d783=\", \"20\u00d768=\", \"25\u00d762=\", \"66\u00d743=\"],\n  [\"69\u00d771=\", \"36\u00d785=\", \"83\u00d731=\", \"38\u00d784=\", \"92\u00d756=\"],\n  [\"11\u00d730=\", \"74\u00d765=\", \"52\u00d779=\", \"93\u00d795=\", \"63\u00d797=\"],\n  [\"98\u00d711=\", \"62\u00d771=\", \"48\u00d729=\", \"77\u00d746=\", \"66\u00d792=\"],\n  [\"47\u00d713=\", \"38\u00d758=\", \"29\u00d723=\", \"30\u00d755=\", \"61\u00d731=\"],\n  [\"86\u00d745=\", \"40\u00d755=\", \"89\u00d764=\", \"80\u00d750=\", \"72\u00d768=\"],\n  [\"19\u00d726=\", \"29\u00d748=\", \"62\u00d792=\", \"42\u00d725=\", \"36\u00d783=\"],\n  [\"37\u00d752=\", \"99\u00d728=\", \"78\u00d785=\", \"64\u00d767=\", \"93\u00d780=\"],\n  [\"96\u00d715=\", \"10\u00d726=\", \"57\u00d744=\", \"59\u00d773=\", \"61\u00d733=\"],\n  [\"78\u00d764=\", \"50\u00d755=\", \"91\u00d771=\", \"62\u00d768=\", \"95\u00d738=\"],\n  [\"42\u00d769=\", \"81\u00d769=\", \"96\u00d726=\", \"44\u00d767=\", \"51\u00d722=\"],\n  [\"91\u00d791=\", \"10\u00d726=\", \"82\u00d766=\", \"94\u00d745=\", \"17\u00d783=\"],\n  [\"86\u00d713=\", \"11\u00d736=\", \"62\u00d789=\", \"83\u00d796=\", \"39\u00d756=\"],\n  [\"22\u00d724=\", \"93\u00d745=\", \"22\u00d743=\", \"82\u00d782=\", \"21\u00d711=\"],\n  [\"12\u00d777=\", \"99\u00d756=\", \"43\u00d711=\", \"21\u00d735=\", \"79\u00d782=\"],\n  [\"61\u00d745=\", \"65\u00d726=\", \"45\u00d753=\", \"86\u00d766=\", \"70\u00d753=\"],\n  [\"71\u00d752=\", \"83\u00d715=\", \"34\u00d762=\", \"91\u00d712=\", \"44\u00d796=\"],\n  [\"81\u00d774=\", \"18\u00d712=\", \"100\u00d752=\", \"60\u00d740=\", \"82\u00d780=\"],\n  [\"60\u00d782=\", \"91\u00d747=\", \"54\u00d765=\", \"44\u00d785=\", \"14\u00d790=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 \"NN\u00d7NN=\" multiplication prompts in the single 20x5 table\n# with the new values from the target revision. Writing directly to each\n# cell's Range.Text preserves the cell/paragraph/run formatting (fonts,\n# size, alignment) and only swaps the visible text.\n$newValues = @(\n    @('97\u00d724=','27\u00d774=','34\u00d785=','61\u00d735=','49\u00d721='),\n    @('54\u00d758=','43\u00d783=','20\u00d768=','25\u00d762=','66\u00d743='),\n    @('69\u00d771=','36\u00d785=','83\u00d731=','38\u00d784=','92\u00d756='),\n    @('11\u00d730=','74\u00d765=','52\u00d779=','93\u00d795=','63\u00d797='),\n    @('98\u00d711=','62\u00d771=','48\u00d729=','77\u00d746=','66\u00d792='),\n    @('47\u00d713=','38\u00d758=','29\u00d723=','30\u00d755=','61\u00d731='),\n    @('86\u00d745=','40\u00d755=','89\u00d764=','80\u00d750=','72\u00d768='),\n    @('19\u00d726=','29\u00d748=','62\u00d792=','42\u00d725=','36\u00d783='),\n    @('37\u00d752=','99\u00d728=','78\u00d785=','64\u00d767=','93\u00d780='),\n    @('96\u00d715=','10\u00d726=','57\u00d744=','59\u00d773=','61\u00d733='),\n    @('78\u00d764=','50\u00d755=','91\u00d771=','62\u00d768=','95\u00d738='),\n    @('42\u00d769=','81\u00d769=','96\u00d726=','44\u00d767=','51\u00d722='),\n    @('91\u00d791=','10\u00d726=','82\u00d766=','94\u00d745=','17\u00d783='),\n    @('86\u00d713=','11\u00d736=','62\u00d789=','83\u00d796=','39\u00d756='),\n    @('22\u00d724=','93\u00d745=','22\u00d743=','82\u00d782=','21\u00d711='),\n    @('12\u00d777=','99\u00d756=','43\u00d711=','21\u00d735=','79\u00d782='),\n    @('61\u00d745=','65\u00d726=','45\u00d753=','86\u00d766=','70\u00d753='),\n    @('71\u00d752=','83\u00d715=','34\u00d762=','91\u00d712=','44\u00d796='),\n    @('81\u00d774=','18\u00d712=','100\u00d752=','60\u00d740=','82\u00d780='),\n    @('60\u00d782=','91\u00d747=','54\u00d765=','44\u00d785=','14\u00d790=')\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
